$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 50
$ws1.Range("F11").Value = 1401
$ws1.Range("F12").Value = 3043
$ws1.Range("F13").Value = 543
$ws1.Range("F16").Value = 826
$ws1.Range("F17").Value = 257
$ws1.Range("F25").Value = 3930
$ws1.Range("F26").Value = 720
$ws1.Range("F30").Value = 68

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 41

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 50
$ws4.Range("F14").Value = 41
$ws4.Range("F22").Value = 1401
$ws4.Range("F24").Value = 543
$ws4.Range("F27").Value = 826
$ws4.Range("F28").Value = 257
$ws4.Range("F38").Value = 3930
$ws4.Range("F39").Value = 720
$ws4.Range("F45").Value = 68
